$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated correlation values for the "sp500" row (row 2) and column (column B),
# reflecting the re-run of the correlation matrix with the updated learning model.

$ws.Range("C2").Value = -0.7515833566942681
$ws.Range("D2").Value = -0.1211002744327765
$ws.Range("E2").Value = -0.03198733403962323
$ws.Range("G2").Value = -0.03067636664815504
$ws.Range("H2").Value = -0.2128317531051811
$ws.Range("I2").Value = -0.04909269184698396
$ws.Range("J2").Value = -0.109523284579041
$ws.Range("K2").Value = 0.07145745119241201
$ws.Range("L2").Value = -0.09569154803023676
$ws.Range("M2").Value = 0.03524498611133097
$ws.Range("N2").Value = -0.238878242919409
$ws.Range("O2").Value = 0.05244044240843464

$ws.Range("B3").Value = -0.7515833566942681
$ws.Range("B4").Value = -0.1211002744327765
$ws.Range("B5").Value = -0.03198733403962323
$ws.Range("B7").Value = -0.03067636664815504
$ws.Range("B8").Value = -0.2128317531051811
$ws.Range("B9").Value = -0.04909269184698396
$ws.Range("B10").Value = -0.109523284579041
$ws.Range("B11").Value = 0.07145745119241201
$ws.Range("B12").Value = -0.09569154803023676
$ws.Range("B13").Value = 0.03524498611133097
$ws.Range("B14").Value = -0.238878242919409
$ws.Range("B15").Value = 0.05244044240843464
